$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 319; this shifts rows 319:382 down to 320:383,
# matching the rest of the worksheet and extending the used range to R383.
$ws.Rows.Item(319).Insert()

# Populate the newly inserted row 319 with its data (same static columns as
# the surrounding "Apio" records for "Macroferia Regional de Talca").
$ws.Range("A319").Value2 = 5
$ws.Range("B319").Value2 = "Macroferia Regional de Talca"
$ws.Range("C319").Value2 = "Maule"
$ws.Range("D319").Value2 = 45211
$ws.Range("E319").Value2 = 7
$ws.Range("F319").Value2 = 100112017
$ws.Range("G319").Value2 = "Apio"
$ws.Range("H319").Value2 = "Americana (o)"
$ws.Range("I319").Value2 = "Primera"
$ws.Range("J319").Value2 = 500
$ws.Range("K319").Value2 = 7000
$ws.Range("L319").Value2 = 7000
$ws.Range("M319").Value2 = 7000
$ws.Range("N319").Value2 = "`$/docena de matas"
$ws.Range("O319").Value2 = "Provincia del Elqu$([char]0xED)"
$ws.Range("P319").Value2 = 1167
$ws.Range("Q319").Value2 = 6
$ws.Range("R319").Value2 = "Hortaliza"
